# Generate Report for Handoff
# Rotates the handoff file id/hash from db321889-0546-4e1d-ae2a-50f3011f14f4
# to 4a63b18c-daab-41cd-84bc-2f45df94d6b0 (and the associated xlf content
# hash from bec8a18adde371b97a065b27a398402bd8ea7e85 to
# 518e73927a685d347fda40c1c9aea2daa2db1f64), and bumps the handoff
# timestamps that were recorded alongside the new files.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldId   = "db321889-0546-4e1d-ae2a-50f3011f14f4"
$newId   = "4a63b18c-daab-41cd-84bc-2f45df94d6b0"

$newMdName    = "$newId.md"
$newMdDisplay = "e2e\$newId.md"

$newZhXlf = "$newId.518e73927a685d347fda40c1c9aea2daa2db1f64.zh-cn.xlf"
$newDeXlf = "$newId.518e73927a685d347fda40c1c9aea2daa2db1f64.de-de.xlf"

$newOverviewDate = "2016-08-17 12:55:17"
$newZhCnDate     = "2016-08-17 12:55:11"

# --- Overview sheet -------------------------------------------------
# A2: bare file name
$wsOverview.Range("A2").Value = $newMdName

# B2: path+name, also the display text of the hyperlink anchored there
$wsOverview.Range("B2").Value = $newMdDisplay
$hOverview = $wsOverview.Range("B2").Hyperlinks.Item(1)
$hOverview.TextToDisplay = $newMdDisplay

# G2: "Latest HO Xliff Generate Date"
$wsOverview.Range("G2").Value = $newOverviewDate

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn.Range("A2").Value = $newMdName
$hZhCn = $wsZhCn.Range("A2").Hyperlinks.Item(1)
$hZhCn.TextToDisplay = $newMdName

# G2: Latest Handoff File (zh-cn xliff)
$wsZhCn.Range("G2").Value = $newZhXlf
# H2: Latest Handoff Datetime
$wsZhCn.Range("H2").Value = $newZhCnDate

# --- de-de sheet ------------------------------------------------------
$wsDeDe.Range("A2").Value = $newMdName
$hDeDe = $wsDeDe.Range("A2").Hyperlinks.Item(1)
$hDeDe.TextToDisplay = $newMdName

# G2: Latest Handoff File (de-de xliff)
$wsDeDe.Range("G2").Value = $newDeXlf
# H2: Latest Handoff Datetime (shares the Overview sheet's timestamp)
$wsDeDe.Range("H2").Value = $newOverviewDate
